$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 30, shifting existing rows (and their formatting) down
$ws.Rows("30:30").Insert()

# Copy full formatting (styles, row height) from the row below (old row 30, now row 31)
$ws.Range("A31:Q31").Copy()
$ws.Range("A30:Q30").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows("30:30").RowHeight = $ws.Rows("31:31").RowHeight

Write-Output "done insert + format copy"
